# Updates to circular messages.
# The "Rectangle 2" circular prompt on the "Trial block completed." slide
# currently reads " Press " + "a button to continue" + "." (three runs).
# Replace it with a single run reading "Doff HMD and press 'Return'."
# (curly quotes), keeping the formatting of the first original run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item("Rectangle 2")

$lq = [char]0x2018
$rq = [char]0x2019
$newText = "Doff HMD and press " + $lq + "Return" + $rq + "."

$tr = $sh.TextFrame.TextRange

# First collapse to an unrelated placeholder string so PowerPoint merges
# the three existing runs into a single run (inheriting the first run's
# formatting) instead of preserving a leftover run for the trailing ".".
$tr.Text = "x"

# Now set the real replacement text on the (now single-run) text range.
$tr2 = $sh.TextFrame.TextRange
$tr2.Text = $newText
